# Fruta / hortaliza, semanal
# Insert two new weekly price rows for Durazno (Vega Modelo de Temuco) at rows 260-261,
# pushing the existing rows 260-328 down to 262-330.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 260 (formatting carried down from row above).
$ws.Range("A260:A261").EntireRow.Insert()

# Row 260: new "Early Majestic" / "Primera" entry
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 44924
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100103
$ws.Cells.Item(260, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(260, 9).Value = 100103004
$ws.Cells.Item(260, 10).Value = "Durazno"
$ws.Cells.Item(260, 11).Value = "Early Majestic"
$ws.Cells.Item(260, 12).Value = "Primera"
$ws.Cells.Item(260, 13).Value = 280
$ws.Cells.Item(260, 14).Value = 20000
$ws.Cells.Item(260, 15).Value = 21000
$ws.Cells.Item(260, 16).Value = 20446
$ws.Cells.Item(260, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(260, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(260, 19).Value = 1136
$ws.Cells.Item(260, 20).Value = 18

# Row 261: new "Early Majestic" / "Segunda" entry
$ws.Cells.Item(261, 1).Value = 10
$ws.Cells.Item(261, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(261, 3).Value = "La Araucanía"
$ws.Cells.Item(261, 4).Value = 44924
$ws.Cells.Item(261, 5).Value = 9
$ws.Cells.Item(261, 6).Value = "Fruta"
$ws.Cells.Item(261, 7).Value = 100103
$ws.Cells.Item(261, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(261, 9).Value = 100103004
$ws.Cells.Item(261, 10).Value = "Durazno"
$ws.Cells.Item(261, 11).Value = "Early Majestic"
$ws.Cells.Item(261, 12).Value = "Segunda"
$ws.Cells.Item(261, 13).Value = 120
$ws.Cells.Item(261, 14).Value = 16000
$ws.Cells.Item(261, 15).Value = 16000
$ws.Cells.Item(261, 16).Value = 16000
$ws.Cells.Item(261, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(261, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(261, 19).Value = 889
$ws.Cells.Item(261, 20).Value = 18
